$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "pitt_louisville_01112025" entry (row 2), shifting the rest of the
# table up by one row.
$ws.Rows.Item(2).Delete()

# Update the selected cell to match the new active selection.
$ws.Range("G7").Select()
